# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet named "2022-Q3" right after "总计" (so the tab
#    order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2020-Q4).
# 2. Populate it with the new quarter's fund-holding table.
# 3. Prepend the new quarter's summary row to the "总计" sheet (shifting
#    the existing rows down by one) and append the 2020-Q4 row that falls
#    out the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a numeric-looking value into a cell while keeping it
# text-typed (so "1.11" stays the string "1.11", not the number 1.11).
# A formula that evaluates to a string literal, copied back in as a
# value, gives Excel a genuinely text-typed result with no left-over
# "quote prefix" cell format.
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$scratch = $summarySheet.Cells.Item(200, 26)   # far-away throwaway cell (Z200)

function Set-TextValue($cell, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet.
#    Copy "2022-Q1" as a template: it already has the right shape (one
#    header row + two data rows) and styling for a quarterly holdings
#    sheet, so the new sheet inherits identical formatting.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q1")
$template.Copy($null, $summarySheet)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Fill in the 2022-Q3 holdings data (overwrite the template's values).
# ---------------------------------------------------------------------
$q3.Cells.Item(2, 1).Value = 0
Set-TextValue $q3.Cells.Item(2, 2) "003359"
Set-TextValue $q3.Cells.Item(2, 3) "大成中证360互联网+大数据100指数C"
Set-TextValue $q3.Cells.Item(2, 4) "1.11"
Set-TextValue $q3.Cells.Item(2, 5) "92.17"
Set-TextValue $q3.Cells.Item(2, 6) "1.00"
Set-TextValue $q3.Cells.Item(2, 7) "0.0111"
$q3.Cells.Item(2, 8).Value = 5

$q3.Cells.Item(3, 1).Value = 1
Set-TextValue $q3.Cells.Item(3, 2) "002236"
Set-TextValue $q3.Cells.Item(3, 3) "大成中证360互联网+大数据100指数A"
Set-TextValue $q3.Cells.Item(3, 4) "1.03"
Set-TextValue $q3.Cells.Item(3, 5) "92.17"
Set-TextValue $q3.Cells.Item(3, 6) "1.00"
Set-TextValue $q3.Cells.Item(3, 7) "0.0103"
$q3.Cells.Item(3, 8).Value = 5

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: insert the 2022-Q3 row at the top
#    of the data, push everything else down one row, and add the
#    2020-Q4 row that now lands on row 6.
# ---------------------------------------------------------------------

# First, grow row 6 with the same look as the rest of column A (index
# style) by copying the format down from row 5, the previous last row.
$summarySheet.Cells.Item(5, 1).Copy()
$summarySheet.Cells.Item(6, 1).PasteSpecial(-4122)   # xlPasteFormats
$summarySheet.Cells.Item(6, 1).Value = 4
Set-TextValue $summarySheet.Cells.Item(6, 2) "2020-Q4"
$summarySheet.Cells.Item(6, 3).Value = 2
$summarySheet.Cells.Item(6, 4).Value = 0

$summarySheet.Cells.Item(5, 1).Value = 3
Set-TextValue $summarySheet.Cells.Item(5, 2) "2021-Q4"
$summarySheet.Cells.Item(5, 3).Value = 5
$summarySheet.Cells.Item(5, 4).Value = 0.99

$summarySheet.Cells.Item(4, 1).Value = 2
Set-TextValue $summarySheet.Cells.Item(4, 2) "2022-Q1"
$summarySheet.Cells.Item(4, 3).Value = 2
$summarySheet.Cells.Item(4, 4).Value = 0.78

$summarySheet.Cells.Item(3, 1).Value = 1
Set-TextValue $summarySheet.Cells.Item(3, 2) "2022-Q2"
$summarySheet.Cells.Item(3, 3).Value = 5
$summarySheet.Cells.Item(3, 4).Value = 0.87

$summarySheet.Cells.Item(2, 1).Value = 0
Set-TextValue $summarySheet.Cells.Item(2, 2) "2022-Q3"
$summarySheet.Cells.Item(2, 3).Value = 2
$summarySheet.Cells.Item(2, 4).Value = 0.02

# ---------------------------------------------------------------------
# Clean up the scratch cell so it leaves no trace in the saved sheet.
# ---------------------------------------------------------------------
$scratch.Clear()

# Copying a sheet shifts Excel's active-tab focus onto the new sheet;
# restore the original active sheet ("2020-Q4", the last tab) so the
# workbook's active-tab/selection state is unchanged by this edit.
$wb.Worksheets.Item("2020-Q4").Activate()
